# "New Features implementation plan-2020.xlsx" update
#  - News story progress bar feature: time estimate revised, and a note
#    added to "Ask Jesse for advice" (kept struck-through like the rest
#    of the original note run).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3: revised time estimate
$ws.Range("D3").Value2 = "1hr 15min"

# E3: insert "Ask Jesse for advice;" into the notes, right after the
# existing (struck-through) "Research mobile applications;" remark.
$ws.Range("E3").Value2 = "Research mobile applications; Ask Jesse for advice; Determine PRL-appropriate structure; Design XD mockup; Maybe ask Federica for opinion"
$ws.Range("E3").Characters(1, 51).Font.Strikethrough = $true
$ws.Range("E3").Characters(52, 86).Font.Strikethrough = $false

# Row 3 grows by one wrapped line to fit the longer note.
$ws.Rows.Item(3).RowHeight = 58

# Leave the cursor where the edit was made.
$ws.Range("B3:C3").Select()
